# Error Calculations and Plots
# This dataset is a "missing data" simulation sheet (combination_3_ABCDF / BC / seed1).
# The edit:
#  1) Drops two rows that were removed from the sample ("RM 232" and "SC 92"),
#     shrinking the sheet from A1:F35 to A1:F33.
#  2) Re-draws which cells in columns C/D count as "missing" (rendered as blank/
#     empty-text cells) for the remaining rows, restoring some values and
#     blanking others.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Remove the two dropped rows -----------------------------------------
# Row 26 = "RM 232" -> delete entirely; everything below shifts up by one.
$ws.Rows("26:26").Delete()
# The old "SC 92" row was at 28, and is now at 27 after the shift above.
$ws.Rows("27:27").Delete()

# --- helper: blank a cell out the same way the original "missing" cells are
# represented (an empty, text-typed cell) instead of clearing it outright.
# A leading "'" forces Excel to store an empty *text* value rather than
# clearing the cell to a true blank; reset the style afterwards so it does
# not leave a stray quote-prefix format behind.
function Set-Missing($addr) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}

# --- 2) Column D re-impute/re-blank for rows 6,8,18,20,23,25 ----------------
$ws.Range("D6").Value = -14.2
Set-Missing "D8"
$ws.Range("D18").Value = -15.2
Set-Missing "D20"
$ws.Range("D23").Value = -13.9
Set-Missing "D25"

# --- 3) Column C/D fixups for the rows following the deletions --------------
# (row numbers below are the FINAL positions, after the two row deletions)
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = -14.6

Set-Missing "C28"
$ws.Range("D28").Value = -13.7

Set-Missing "C29"
$ws.Range("D29").Value = -13

$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6

Set-Missing "C32"
$ws.Range("D32").Value = -14.7
